$wb = $excel.ActiveWorkbook

# --- Sheet "spec": fix a few risk-matrix cells (value + matching fill color) ---
$spec = $wb.Worksheets.Item("spec")
$spec.Activate()

# H2: 1 -> 2 (recolor to the "2" level color, taken from a stable s=13 cell)
$spec.Range("H2").Interior.Color = $spec.Range("I2").Interior.Color
$spec.Range("H2").Value = 2

# G3: 0 -> 1 (recolor to the "1" level color, taken from a stable s=14 cell)
$spec.Range("G3").Interior.Color = $spec.Range("H3").Interior.Color
$spec.Range("G3").Value = 1

# J4: 2 -> 1 (recolor to the "1" level color, taken from a stable s=14 cell)
$spec.Range("J4").Interior.Color = $spec.Range("I4").Interior.Color
$spec.Range("J4").Value = 1

# I5: 1 -> 0 (recolor to the "0" level color, taken from a stable s=15 cell)
$spec.Range("I5").Interior.Color = $spec.Range("H5").Interior.Color
$spec.Range("I5").Value = 0

# Restore the view state for this sheet (scroll + selection)
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$spec.Range("I5").Select()

# --- Sheet "library_content": bump the severity value and move the selection ---
$lib = $wb.Worksheets.Item("library_content")
$lib.Activate()
$lib.Range("B2").Value = 4
$lib.Range("B7").Select()
